$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 375 (shifts existing rows 375-454 down to 376-455)
$ws.Rows(375).Insert()

# Populate the newly inserted row 375 with its data
$ws.Cells.Item(375, 1).Value  = 10
$ws.Cells.Item(375, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(375, 3).Value  = "La Araucanía"
$ws.Cells.Item(375, 4).Value  = 44995
$ws.Cells.Item(375, 5).Value  = 9
$ws.Cells.Item(375, 6).Value  = 100112044
$ws.Cells.Item(375, 7).Value  = "Perejil"
$ws.Cells.Item(375, 8).Value  = "Sin especificar"
$ws.Cells.Item(375, 9).Value  = "Primera"
$ws.Cells.Item(375, 10).Value = 45
$ws.Cells.Item(375, 11).Value = 5000
$ws.Cells.Item(375, 12).Value = 5000
$ws.Cells.Item(375, 13).Value = 5000
$ws.Cells.Item(375, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(375, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(375, 16).Value = 1667
$ws.Cells.Item(375, 17).Value = 3
$ws.Cells.Item(375, 18).Value = "Hortaliza"
